$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44289.25
$ws.Range("J3").Value = 44289.25
$ws.Range("L3").Value = 44289.25
$ws.Range("N3").Value = -44517.25

$ws.Range("H64").Value = 13125
$ws.Range("J64").Value = 13125
$ws.Range("L64").Value = 13125
$ws.Range("N64").Value = -13621

$ws.Range("H67").Value = 13125
$ws.Range("J67").Value = 13125
$ws.Range("L67").Value = 13125
$ws.Range("N67").Value = -14841

$ws.Range("H102").Value = 44289.25
$ws.Range("J102").Value = 44289.25
$ws.Range("L102").Value = 44289.25
$ws.Range("N102").Value = -50779.25

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws.Range("H141").Value = 5500
$ws.Range("J141").Value = 5500
$ws.Range("L141").Value = 16500
$ws.Range("N141").Value = -26860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 933.3333

$ws.Range("H97").Value = 2454.75
$ws.Range("I97").Value = 2454.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2454.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1958.75
$ws.Range("N97").ClearContents()

$ws.Range("H106").Value = 27997.5
$ws.Range("J106").Value = 27997.5
$ws.Range("L106").Value = 27997.5
$ws.Range("N106").Value = -30521.5

$ws.Range("H116").Value = 933.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 933.3333

$ws.Range("H107").Value = 2492
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 3035
$ws.Range("I134").Value = 3035
$ws.Range("K134").Value = 9105
$ws.Range("M134").Value = -6570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 52000
$ws.Range("J18").Value = 52000
$ws.Range("L18").Value = 52000
$ws.Range("N18").Value = -52460

$ws.Range("H25").Value = 1724.75
$ws.Range("I25").Value = 1724.75
$ws.Range("K25").Value = 1724.75
$ws.Range("M25").Value = -1550.75

$ws.Range("H39").Value = 7567.1113
$ws.Range("J39").Value = 14885.667
$ws.Range("L39").Value = 14885.667
$ws.Range("N39").Value = -15667.667

$ws.Range("H42").Value = 26666
$ws.Range("J42").Value = 26666
$ws.Range("L42").Value = 26666
$ws.Range("N42").Value = -27852

$ws.Range("H49").Value = 7567.1113
$ws.Range("J49").Value = 14885.667
$ws.Range("L49").Value = 14885.667
$ws.Range("N49").Value = -15249.667

$ws.Range("H50").Value = 49196.75
$ws.Range("J50").Value = 49196.75
$ws.Range("L50").Value = 49196.75
$ws.Range("N50").Value = -50446.75

$ws.Range("H51").Value = 45000
$ws.Range("J51").Value = 45000
$ws.Range("L51").Value = 45000
$ws.Range("N51").Value = -46472

$ws.Range("H60").Value = 20052.777
$ws.Range("I60").Value = 7447.5
$ws.Range("K60").Value = 7447.5
$ws.Range("M60").Value = -6936.5

$ws.Range("H61").Value = 45000
$ws.Range("J61").Value = 45000
$ws.Range("L61").Value = 45000
$ws.Range("N61").Value = -45696

$ws.Range("H94").Value = 850
$ws.Range("I94").Value = 850
$ws.Range("K94").Value = 850
$ws.Range("M94").Value = -399

$ws.Range("H104").Value = 52000
$ws.Range("J104").Value = 52000
$ws.Range("L104").Value = 52000
$ws.Range("N104").Value = -57242

$ws.Range("H106").Value = 40791.75
$ws.Range("J106").Value = 40791.75
$ws.Range("L106").Value = 40791.75
$ws.Range("N106").Value = -43315.75

$ws.Range("H109").Value = 69999
$ws.Range("J109").Value = 69999
$ws.Range("L109").Value = 69999
$ws.Range("N109").Value = -72079

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388

$ws.Range("H18").Value = 25000
$ws.Range("J18").Value = 25000
$ws.Range("L18").Value = 25000
$ws.Range("N18").Value = -25586

$ws.Range("H99").Value = 19000
$ws.Range("I99").Value = 13500
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 13500
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -11254
$ws.Range("N99").Value = -34492

$ws.Range("H102").Value = 52437.375
$ws.Range("I102").Value = 52437.375
$ws.Range("K102").Value = 52437.375
$ws.Range("M102").Value = -50815.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -392

$ws.Range("H46").Value = 287070.28
$ws.Range("I46").Value = 500750
$ws.Range("J46").Value = 2164
$ws.Range("K46").Value = 500750
$ws.Range("L46").Value = 2164
$ws.Range("M46").Value = -500562
$ws.Range("N46").Value = -2540

$ws.Range("H116").Value = 68000
$ws.Range("J116").Value = 68000
$ws.Range("L116").Value = 68000
$ws.Range("N116").Value = -77178

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 29996
$ws.Range("I51").Value = 29989
$ws.Range("J51").Value = 29999.5
$ws.Range("K51").Value = 29989
$ws.Range("L51").Value = 29999.5
$ws.Range("M51").Value = -29479
$ws.Range("N51").Value = -31019.5

$ws.Range("H52").Value = 39495
$ws.Range("I52").Value = 38990
$ws.Range("J52").Value = 40000
$ws.Range("K52").Value = 38990
$ws.Range("L52").Value = 40000
$ws.Range("M52").Value = -38764
$ws.Range("N52").Value = -40452

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H81").Value = 855.3333
$ws.Range("I81").Value = 855.3333
$ws.Range("K81").Value = 1710.6666
$ws.Range("M81").Value = -649.6666

$ws.Range("H84").Value = 855.3333
$ws.Range("I84").Value = 855.3333
$ws.Range("K84").Value = 8553.333000000001
$ws.Range("M84").Value = -3249.333000000001

$ws.Range("H96").Value = 42166.5
$ws.Range("I96").Value = 14000
$ws.Range("K96").Value = 14000
$ws.Range("M96").Value = -12627

$ws.Range("H97").Value = 61000
$ws.Range("J97").Value = 61000
$ws.Range("L97").Value = 61000
$ws.Range("N97").Value = -62982

$ws.Range("H109").Value = 53999.5
$ws.Range("J109").Value = 53999.5
$ws.Range("L109").Value = 53999.5
$ws.Range("N109").Value = -56773.5
